# Auto update stock data
#
# Updates the "Date_1" and "EBITDA" columns (A and B) for the first data
# row of each of the 13 companies in the sheet (rows 2,8,14,...,74 - every
# 6th row), moving the snapshot date from 2025/11/24 to 2025/11/25 and
# refreshing the EBITDA figure. All of these sheet cells are stored as
# plain text (not real dates/numbers), so we briefly force a Text number
# format while writing the value (to stop Excel from auto-converting the
# string into a date serial / numeric value) and then restore the cell's
# style back to Normal so no residual formatting is left on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param([string]$Address, [string]$Text)
    $cell = $ws.Range($Address)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.Style = "Normal"
}

$updates = @(
    @{ Row = 2;  Date = "2025/11/25"; Ebitda = "4.72" },
    @{ Row = 8;  Date = "2025/11/25"; Ebitda = "7.48" },
    @{ Row = 14; Date = "2025/11/25"; Ebitda = "2.81" },
    @{ Row = 20; Date = "2025/11/25"; Ebitda = "12.13" },
    @{ Row = 26; Date = "2025/11/25"; Ebitda = "9.79" },
    @{ Row = 32; Date = "2025/11/25"; Ebitda = $null },
    @{ Row = 38; Date = "2025/11/25"; Ebitda = "38.59" },
    @{ Row = 44; Date = "2025/11/25"; Ebitda = "10.20" },
    @{ Row = 50; Date = "2025/11/25"; Ebitda = "11.41" },
    @{ Row = 56; Date = "2025/11/25"; Ebitda = "33.66" },
    @{ Row = 62; Date = "2025/11/25"; Ebitda = "10.77" },
    @{ Row = 68; Date = "2025/11/25"; Ebitda = "11.95" },
    @{ Row = 74; Date = "2025/11/25"; Ebitda = "15.58" }
)

foreach ($u in $updates) {
    $dateAddr = "A" + $u.Row
    Set-TextValue $dateAddr $u.Date
    if ($null -ne $u.Ebitda) {
        $ebitdaAddr = "B" + $u.Row
        Set-TextValue $ebitdaAddr $u.Ebitda
    }
}
